# Update profit files after running on 2026-02-15
# Appends one new data row (row 83) to Sheet1, mirroring the existing
# daily-snapshot rows above it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 83

# Column A holds dates stored as literal text (e.g. "02/14/2026" in A82),
# not real Excel dates. A leading apostrophe forces the new value to be
# entered as text too, instead of being auto-converted to a date serial.
$ws.Cells.Item($row, 1).Value = "'02/15/2026"

$ws.Cells.Item($row, 2).Value  = 9598.16
$ws.Cells.Item($row, 3).Value  = 0.242646236788244
$ws.Cells.Item($row, 4).Value  = 0.757353763211756
$ws.Cells.Item($row, 5).Value  = -307.72
$ws.Cells.Item($row, 6).Value  = -34.77
$ws.Cells.Item($row, 7).Value  = -23672.56
$ws.Cells.Item($row, 8).Value  = -76.51000000000001
$ws.Cells.Item($row, 9).Value  = -1073.89
$ws.Cells.Item($row, 10).Value = -31.56
$ws.Cells.Item($row, 11).Value = -24746.45
$ws.Cells.Item($row, 12).Value = -72.05
